# Apply edits described in the commit "add and more results"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update timing values in rows 5-7 (column B)
$ws.Cells.Item(5, 2).Value = 0.0005500316619873047
$ws.Cells.Item(6, 2).Value = 0.0004568099975585938
$ws.Cells.Item(7, 2).Value = 0.002872228622436523

# 2. Switch the embedding tuples in column A (rows 8, 49, 104) from
#    Python-tuple-style "(a, b)" text to list-style "[a, b]" text.
$ws.Cells.Item(8, 1).Value = "[[2, 0], [2, 2], [1, 3], [0, 3], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Cells.Item(49, 1).Value = "[[2, 0], [2, 1], [1, 1], [1, 2], [0, 2], [0, 3], [1, 3], [0, 1], [2, 2], [1, 0], [3, 1]]"
$ws.Cells.Item(104, 1).Value = "[[1, 1], [1, 2], [2, 1], [1, 0], [0, 1], [2, 2], [0, 2], [0, 0], [2, 0], [0, 3], [1, 3]]"

# 3. Insert a new row at 154 (this pushes the former rows 154-158 down to 155-159,
#    and Excel auto-extends the used range/dimension to N159).
$ws.Rows.Item(154).Insert()

# 4. Populate the newly inserted row 154 with the "move_fidelity" result.
$ws.Cells.Item(154, 1).Value = "move_fidelity"
$ws.Cells.Item(154, 2).Value = 0.9990331134467497

# 5. Update the "total time:" row (now row 158) with its new value.
$ws.Cells.Item(158, 2).Value = 0.01346921920776367
